$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Rows.Item(1).Cells.Item(1).Range.Text = "45+39="
$t.Rows.Item(1).Cells.Item(2).Range.Text = "91-46="
$t.Rows.Item(1).Cells.Item(3).Range.Text = "45+16="
$t.Rows.Item(1).Cells.Item(4).Range.Text = "19+4="
$t.Rows.Item(1).Cells.Item(5).Range.Text = "91-73="
$t.Rows.Item(2).Cells.Item(1).Range.Text = "5+87="
$t.Rows.Item(2).Cells.Item(2).Range.Text = "93-75="
$t.Rows.Item(2).Cells.Item(3).Range.Text = "63-17="
$t.Rows.Item(2).Cells.Item(4).Range.Text = "89+4="
$t.Rows.Item(2).Cells.Item(5).Range.Text = "82-78="
$t.Rows.Item(3).Cells.Item(1).Range.Text = "7+64="
$t.Rows.Item(3).Cells.Item(2).Range.Text = "63-17="
$t.Rows.Item(3).Cells.Item(3).Range.Text = "58+34="
$t.Rows.Item(3).Cells.Item(4).Range.Text = "70-34="
$t.Rows.Item(3).Cells.Item(5).Range.Text = "93-48="
$t.Rows.Item(4).Cells.Item(1).Range.Text = "29+12="
$t.Rows.Item(4).Cells.Item(2).Range.Text = "86-78="
$t.Rows.Item(4).Cells.Item(3).Range.Text = "9+45="
$t.Rows.Item(4).Cells.Item(4).Range.Text = "46+25="
$t.Rows.Item(4).Cells.Item(5).Range.Text = "92-64="
$t.Rows.Item(5).Cells.Item(1).Range.Text = "4+79="
$t.Rows.Item(5).Cells.Item(2).Range.Text = "92-56="
$t.Rows.Item(5).Cells.Item(3).Range.Text = "28+44="
$t.Rows.Item(5).Cells.Item(4).Range.Text = "38+8="
$t.Rows.Item(5).Cells.Item(5).Range.Text = "73-24="
$t.Rows.Item(6).Cells.Item(1).Range.Text = "54+38="
$t.Rows.Item(6).Cells.Item(2).Range.Text = "41-13="
$t.Rows.Item(6).Cells.Item(3).Range.Text = "67-39="
$t.Rows.Item(6).Cells.Item(4).Range.Text = "50-17="
$t.Rows.Item(6).Cells.Item(5).Range.Text = "9+37="
$t.Rows.Item(7).Cells.Item(1).Range.Text = "94-25="
$t.Rows.Item(7).Cells.Item(2).Range.Text = "35-16="
$t.Rows.Item(7).Cells.Item(3).Range.Text = "32-7="
$t.Rows.Item(7).Cells.Item(4).Range.Text = "7+47="
$t.Rows.Item(7).Cells.Item(5).Range.Text = "47+6="
$t.Rows.Item(8).Cells.Item(1).Range.Text = "70-55="
$t.Rows.Item(8).Cells.Item(2).Range.Text = "17+17="
$t.Rows.Item(8).Cells.Item(3).Range.Text = "48+47="
$t.Rows.Item(8).Cells.Item(4).Range.Text = "77-9="
$t.Rows.Item(8).Cells.Item(5).Range.Text = "93-77="
$t.Rows.Item(9).Cells.Item(1).Range.Text = "23+59="
$t.Rows.Item(9).Cells.Item(2).Range.Text = "61-46="
$t.Rows.Item(9).Cells.Item(3).Range.Text = "82-44="
$t.Rows.Item(9).Cells.Item(4).Range.Text = "3+68="
$t.Rows.Item(9).Cells.Item(5).Range.Text = "93-36="
$t.Rows.Item(10).Cells.Item(1).Range.Text = "21-6="
$t.Rows.Item(10).Cells.Item(2).Range.Text = "85-57="
$t.Rows.Item(10).Cells.Item(3).Range.Text = "18+53="
$t.Rows.Item(10).Cells.Item(4).Range.Text = "52+39="
$t.Rows.Item(10).Cells.Item(5).Range.Text = "7+76="
$t.Rows.Item(11).Cells.Item(1).Range.Text = "87+8="
$t.Rows.Item(11).Cells.Item(2).Range.Text = "50-8="
$t.Rows.Item(11).Cells.Item(3).Range.Text = "50-21="
$t.Rows.Item(11).Cells.Item(4).Range.Text = "75-19="
$t.Rows.Item(11).Cells.Item(5).Range.Text = "9+86="
$t.Rows.Item(12).Cells.Item(1).Range.Text = "91-69="
$t.Rows.Item(12).Cells.Item(2).Range.Text = "54+29="
$t.Rows.Item(12).Cells.Item(3).Range.Text = "77-18="
$t.Rows.Item(12).Cells.Item(4).Range.Text = "80-43="
$t.Rows.Item(12).Cells.Item(5).Range.Text = "48+14="
$t.Rows.Item(13).Cells.Item(1).Range.Text = "90-28="
$t.Rows.Item(13).Cells.Item(2).Range.Text = "97-79="
$t.Rows.Item(13).Cells.Item(3).Range.Text = "82-29="
$t.Rows.Item(13).Cells.Item(4).Range.Text = "90-6="
$t.Rows.Item(13).Cells.Item(5).Range.Text = "36+55="

# Merge old rows 14 and 15 (1-indexed) into a single row 14, then delete row 15
$t.Rows.Item(14).Cells.Item(1).Range.Text = "59+7="
$t.Rows.Item(14).Cells.Item(2).Range.Text = "17+38="
$t.Rows.Item(14).Cells.Item(3).Range.Text = "14+28="
$t.Rows.Item(14).Cells.Item(4).Range.Text = "25-8="
$t.Rows.Item(14).Cells.Item(5).Range.Text = "18+36="
$t.Rows.Item(15).Delete()

# Remaining rows 15-19 (1-indexed, post-delete) map from before rows 15-19 (0-indexed, i.e. old physical 16-20)
$t.Rows.Item(15).Cells.Item(1).Range.Text = "61-54="
$t.Rows.Item(15).Cells.Item(2).Range.Text = "17+46="
$t.Rows.Item(15).Cells.Item(3).Range.Text = "38+44="
$t.Rows.Item(15).Cells.Item(4).Range.Text = "47+36="
$t.Rows.Item(15).Cells.Item(5).Range.Text = "72-53="
$t.Rows.Item(16).Cells.Item(1).Range.Text = "36+45="
$t.Rows.Item(16).Cells.Item(2).Range.Text = "15+29="
$t.Rows.Item(16).Cells.Item(3).Range.Text = "67-19="
$t.Rows.Item(16).Cells.Item(4).Range.Text = "20-18="
$t.Rows.Item(16).Cells.Item(5).Range.Text = "72-5="
$t.Rows.Item(17).Cells.Item(1).Range.Text = "15+46="
$t.Rows.Item(17).Cells.Item(2).Range.Text = "58+39="
$t.Rows.Item(17).Cells.Item(3).Range.Text = "14+8="
$t.Rows.Item(17).Cells.Item(4).Range.Text = "46+48="
$t.Rows.Item(17).Cells.Item(5).Range.Text = "40-15="
$t.Rows.Item(18).Cells.Item(1).Range.Text = "39+44="
$t.Rows.Item(18).Cells.Item(2).Range.Text = "76+15="
$t.Rows.Item(18).Cells.Item(3).Range.Text = "79+5="
$t.Rows.Item(18).Cells.Item(4).Range.Text = "85-49="
$t.Rows.Item(18).Cells.Item(5).Range.Text = "70-68="
$t.Rows.Item(19).Cells.Item(1).Range.Text = "32-25="
$t.Rows.Item(19).Cells.Item(2).Range.Text = "82-34="
$t.Rows.Item(19).Cells.Item(3).Range.Text = "38+16="
$t.Rows.Item(19).Cells.Item(4).Range.Text = "95-88="
$t.Rows.Item(19).Cells.Item(5).Range.Text = "94-89="

# Append a brand-new row 20 with 5 new problems
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "35-19="
$newRow.Cells.Item(2).Range.Text = "28-19="
$newRow.Cells.Item(3).Range.Text = "36+38="
$newRow.Cells.Item(4).Range.Text = "17+6="
$newRow.Cells.Item(5).Range.Text = "82-49="
